$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.987.50"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "'3.467.52"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'591.95"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").Value = "'175.55"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.585"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").Value = "'7.07"
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("D11").Value = "'0.425"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").Value = "'4.067.90"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "'30.62"
$ws.Range("E13").Value = "  +6.50%  "
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "'67.029.97"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("E16").Value = "  -3.82%  "
$ws.Range("D17").Value = "'3.464.52"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "'6.21"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "'14.29"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "'385.89"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("D21").Value = "'7.83"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "'72.45"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'0.532"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("D27").Value = "'10.28"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'6.07"
$ws.Range("E30").Value = "  -3.63%  "
$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("D32").Value = "'2.03"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").Value = "'23.37"
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("D34").Value = "'7.25"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "'162.96"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").Value = "'0.871"
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("D38").Value = "'1.91"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").Value = "'27.19"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").Value = "'4.60"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").Value = "'26.30"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").Value = "'2.784.79"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "'0.0719"
$ws.Range("E44").Value = "  -3.88%  "
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("D46").Value = "'42.13"
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("D47").Value = "'0.0298"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("D48").Value = "'338.39"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").Value = "'1.06"
$ws.Range("E49").Value = "  -3.20%  "
$ws.Range("D50").Value = "'33.10"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").Value = "'6.34"
$ws.Range("E51").Value = "  -2.73%  "
